# Update: Add data for 2021-12-15
# Carjacking arrests by month YoY - advance "through" date from 12-06 to 12-07
# and update the December-to-date and Total rows with the newly incorporated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet / tab name ---------------------------------------------------
$ws.Name = "Through 2021-12-07"

# --- Row 13 (November) ---------------------------------------------------
$ws.Range("U13").Value = 195
$ws.Range("V13").Value = 0.025

# --- Row 14 (December, through date) -------------------------------------
$ws.Range("A14").Value = "December (through 12-07)"

$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 4
# D14 (arrest_rate for 2015) is unchanged

$ws.Range("F14").Value = 21
$ws.Range("G14").Value = 0.0455

$ws.Range("I14").Value = 22
$ws.Range("J14").Value = 0.12

$ws.Range("L14").Value = 13
$ws.Range("M14").Value = 0.07140000000000001

# N14/P14 (2019 columns) did not previously have any data in this row;
# give them the same percentage number format as the other rate columns.
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 7
$ws.Range("P14").Value = 0.125
$ws.Range("P14").NumberFormat = $ws.Range("M14").NumberFormat()

$ws.Range("R14").Value = 33
$ws.Range("S14").Value = 0.0571

$ws.Range("U14").Value = 57

# --- Row 15 (Total) -------------------------------------------------------
$ws.Range("B15").Value = 35
$ws.Range("C15").Value = 262
$ws.Range("D15").Value = 0.1178

$ws.Range("F15").Value = 524
$ws.Range("G15").Value = 0.1043

$ws.Range("I15").Value = 780
$ws.Range("J15").Value = 0.078

$ws.Range("L15").Value = 621
$ws.Range("M15").Value = 0.1078

$ws.Range("N15").Value = 55
$ws.Range("O15").Value = 487
$ws.Range("P15").Value = 0.1015

$ws.Range("R15").Value = 1233
$ws.Range("S15").Value = 0.0508

$ws.Range("U15").Value = 1600
$ws.Range("V15").Value = 0.0583
